$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 13.16594766666667
$ws.Range("H2").Value = 39.497843
$ws.Range("I2").Value = 0.6940777873489595
$ws.Range("J2").Value = 0.6940777873489595
$ws.Range("M2").Value = 0.0006136666666666667
$ws.Range("N2").Value = 0.001841
$ws.Range("O2").Value = 0.000129696697123199
$ws.Range("P2").Value = 0.000129696697123199
$ws.Range("Q2").Value = 0.008079503218111111
$ws.Range("R2").Value = 0.072715528963
$ws.Range("S2").Value = 0.00009001959656573815
$ws.Range("T2").Value = 0.00009001959656573815
$ws.Range("G3").Value = 13.16594766666667
$ws.Range("H3").Value = 39.497843
$ws.Range("I3").Value = 0.6940777873489595
$ws.Range("J3").Value = 0.6940777873489595
$ws.Range("O3").Value = 0.8077214410831794
$ws.Range("P3").Value = 0.8077214410831794
$ws.Range("Q3").Value = 50.31730280972267
$ws.Range("R3").Value = 452.855725287504
$ws.Range("S3").Value = 0.5606215106213261
$ws.Range("T3").Value = 0.5606215106213261
$ws.Range("G4").Value = 13.16594766666667
$ws.Range("H4").Value = 39.497843
$ws.Range("I4").Value = 0.6940777873489595
$ws.Range("J4").Value = 0.6940777873489595
$ws.Range("O4").Value = 0.1921488622196973
$ws.Range("P4").Value = 0.1921488622196973
$ws.Range("Q4").Value = 11.96998370117122
$ws.Range("R4").Value = 107.729853310541
$ws.Range("S4").Value = 0.1333662571310676
$ws.Range("T4").Value = 0.1333662571310676
$ws.Range("I5").Value = 0.1706596770095176
$ws.Range("J5").Value = 0.1706596770095176
$ws.Range("M5").Value = 0.0006136666666666667
$ws.Range("N5").Value = 0.001841
$ws.Range("O5").Value = 0.000129696697123199
$ws.Range("P5").Value = 0.000129696697123199
$ws.Range("Q5").Value = 0.00198658628
$ws.Range("R5").Value = 0.01787927652
$ws.Range("S5").Value = 0.00002213399644024638
$ws.Range("T5").Value = 0.00002213399644024638
$ws.Range("I6").Value = 0.1706596770095176
$ws.Range("J6").Value = 0.1706596770095176
$ws.Range("O6").Value = 0.8077214410831794
$ws.Range("P6").Value = 0.8077214410831794
$ws.Range("S6").Value = 0.1378454802489175
$ws.Range("T6").Value = 0.1378454802489175
$ws.Range("I7").Value = 0.1706596770095176
$ws.Range("J7").Value = 0.1706596770095176
$ws.Range("O7").Value = 0.1921488622196973
$ws.Range("P7").Value = 0.1921488622196973
$ws.Range("S7").Value = 0.03279206276415985
$ws.Range("T7").Value = 0.03279206276415985
$ws.Range("H8").Value = 7.697376999999999
$ws.Range("I8").Value = 0.1352625356415228
$ws.Range("J8").Value = 0.1352625356415228
$ws.Range("M8").Value = 0.0006136666666666667
$ws.Range("N8").Value = 0.001841
$ws.Range("O8").Value = 0.000129696697123199
$ws.Range("P8").Value = 0.000129696697123199
$ws.Range("Q8").Value = 0.001574541228555555
$ws.Range("R8").Value = 0.014170871057
$ws.Range("S8").Value = 0.00001754310411721449
$ws.Range("T8").Value = 0.0000175431041172145
$ws.Range("H9").Value = 7.697376999999999
$ws.Range("I9").Value = 0.1352625356415228
$ws.Range("J9").Value = 0.1352625356415228
$ws.Range("O9").Value = 0.8077214410831794
$ws.Range("P9").Value = 0.8077214410831794
$ws.Range("Q9").Value = 9.805883560517332
$ws.Range("R9").Value = 88.25295204465598
$ws.Range("S9").Value = 0.1092544502129357
$ws.Range("T9").Value = 0.1092544502129357
$ws.Range("H10").Value = 7.697376999999999
$ws.Range("I10").Value = 0.1352625356415228
$ws.Range("J10").Value = 0.1352625356415228
$ws.Range("O10").Value = 0.1921488622196973
$ws.Range("P10").Value = 0.1921488622196973
$ws.Range("S10").Value = 0.02599054232446986
$ws.Range("T10").Value = 0.02599054232446986
